# Update the "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value (same updates apply to both sheets)
$updates = @{
    2  = 8034
    3  = 7683
    9  = 108
    10 = 155
    11 = 224
    13 = 122
    14 = 1252
    16 = 46
    17 = 8
    19 = 104
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
